$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Updates to existing rows 2-23 (ticker symbols changed / cleared / added)
# Each entry only lists the columns that actually change relative to the
# "before" state. $null means "clear the cell".
# ---------------------------------------------------------------------------
$updates = @(
    @{Row=2;  B="NSE:AAVAS";      C="NSE:3IINFOLTD"; E="NSE:BHEL";     F="NSE:ALKEM"},
    @{Row=3;  B="NSE:ADVENZYMES"; C="NSE:BAJFINANCE"; E="NSE:POLYCAB"; F="NSE:IEX"},
    @{Row=4;  B="NSE:ALKEM";      C="NSE:CIEINDIA";   F="NSE:LUPIN"},
    @{Row=5;  B="NSE:ALOKINDS";   C="NSE:CYIENTDLM";  F="NSE:PVRINOX"},
    @{Row=6;  B="NSE:ALPA";       C="NSE:HDFCQUAL";   F="NSE:RAMCOCEM"},
    @{Row=7;  B="NSE:AMIORG";     C="NSE:HILTON"},
    @{Row=8;  B="NSE:ASHOKA";     C="NSE:ITI"},
    @{Row=9;  B="NSE:BLISSGVS";   C="NSE:KKCL"},
    @{Row=10; B="NSE:CAPLIPOINT"; C="NSE:RAILTEL"},
    @{Row=11; B="NSE:CHALET";     C=$null},
    @{Row=12; B="NSE:CLEDUCATE";  C=$null},
    @{Row=13; B="NSE:DEVYANI";    C=$null},
    @{Row=14; B="NSE:GRAVITA";    C=$null},
    @{Row=15; B="NSE:HERITGFOOD"; C=$null},
    @{Row=16; B="NSE:IEX";        C=$null},
    @{Row=17; B="NSE:JBMA";       C=$null},
    @{Row=18; B="NSE:JKCEMENT";   C=$null},
    @{Row=19; B="NSE:JYOTHYLAB";  C=$null},
    @{Row=20; B="NSE:KEC";        C=$null},
    @{Row=21; B="NSE:KEYFINSERV"; C=$null},
    @{Row=22; B="NSE:KOPRAN"},
    @{Row=23; B="NSE:LAL"}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) {
        if ($u.B -eq $null) { $ws.Cells.Item($r, 2).Value = "" } else { $ws.Cells.Item($r, 2).Value = $u.B }
    }
    if ($u.ContainsKey("C")) {
        if ($u.C -eq $null) { $ws.Cells.Item($r, 3).Value = "" } else { $ws.Cells.Item($r, 3).Value = $u.C }
    }
    if ($u.ContainsKey("D")) {
        if ($u.D -eq $null) { $ws.Cells.Item($r, 4).Value = "" } else { $ws.Cells.Item($r, 4).Value = $u.D }
    }
    if ($u.ContainsKey("E")) {
        if ($u.E -eq $null) { $ws.Cells.Item($r, 5).Value = "" } else { $ws.Cells.Item($r, 5).Value = $u.E }
    }
    if ($u.ContainsKey("F")) {
        if ($u.F -eq $null) { $ws.Cells.Item($r, 6).Value = "" } else { $ws.Cells.Item($r, 6).Value = $u.F }
    }
}

# ---------------------------------------------------------------------------
# New rows 24-39, each with a serial number in column A and a ticker in B.
# Column A needs the same style (bold, centered, bordered) as the existing
# numbered rows above it, so copy formatting from row 23's A cell first.
# ---------------------------------------------------------------------------
$newRows = @(
    @{Row=24; A=22; B="NSE:LASA"},
    @{Row=25; A=23; B="NSE:LAURUSLABS"},
    @{Row=26; A=24; B="NSE:LOKESHMACH"},
    @{Row=27; A=25; B="NSE:LUPIN"},
    @{Row=28; A=26; B="NSE:MGL"},
    @{Row=29; A=27; B="NSE:NH"},
    @{Row=30; A=28; B="NSE:NIRAJ"},
    @{Row=31; A=29; B="NSE:NUVOCO"},
    @{Row=32; A=30; B="NSE:PGEL"},
    @{Row=33; A=31; B="NSE:PVRINOX"},
    @{Row=34; A=32; B="NSE:RAMASTEEL"},
    @{Row=35; A=33; B="NSE:RAMCOCEM"},
    @{Row=36; A=34; B="NSE:RANEHOLDIN"},
    @{Row=37; A=35; B="NSE:RBL"},
    @{Row=38; A=36; B="NSE:ROUTE"},
    @{Row=39; A=37; B="NSE:SAGCEM"}
)

$ws.Range("A23").Copy()
foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $nr.A
    $ws.Cells.Item($r, 2).Value = $nr.B
}

$excel.CutCopyMode = 0
